# Update users.xlsx:
#  - rename the worksheet from "Sheet1" to "Users"
#  - replace the two Cyrillic usernames with "user2" / "user3"
#  - leave the rest of the table (passwords, header) untouched
#  - move the active selection to C14, matching the author's last
#    interaction with the sheet before saving

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab) from Sheet1 -> Users
$ws.Name = "Users"

# Update the two username cells that changed
$ws.Range("A3").Value = "user2"
$ws.Range("A4").Value = "user3"

# Move / record the active selection as it was left after editing
[void]$ws.Range("C14").Select()
